$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(45, 1).Value = 44
$ws.Cells.Item(45, 2).Value = ""
$ws.Cells.Item(45, 3).Value = 46
$ws.Cells.Item(45, 4).Value = 18
$ws.Cells.Item(45, 5).Value = "System"
$ws.Cells.Item(45, 6).Value = "2025-03-03 19:39:50"
$ws.Cells.Item(45, 7).Value = 0
